# Update the "as_of_utc" timestamp column (AA) on both data sheets.
$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-12-16 07:07:20"
$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Cells.Item($row, 27).Value = $newTimestamp
    }
}
